# Auto-generated script to update currentAveragePrice / Leve price & profit figures
# across the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 900
$ws.Range("I12").Value = 800
$ws.Range("K12").Value = 800
$ws.Range("M12").Value = -630
$ws.Range("H33").Value = 334.76923
$ws.Range("I33").Value = 304.72726
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 304.72726
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = -75.72726
$ws.Range("N33").Value = -958
$ws.Range("H113").Value = 7110.8887
$ws.Range("I113").Value = 6500
$ws.Range("J113").Value = 7187.25
$ws.Range("K113").Value = 6500
$ws.Range("L113").Value = 7187.25
$ws.Range("M113").Value = -3246
$ws.Range("N113").Value = -13695.25
$ws.Range("H138").Value = 3472.08
$ws.Range("I138").Value = 655.29034
$ws.Range("J138").Value = 4737.594
$ws.Range("K138").Value = 1965.87102
$ws.Range("L138").Value = 14212.782
$ws.Range("M138").Value = 3174.12898
$ws.Range("N138").Value = -24492.782

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1157.85
$ws.Range("I2").Value = 1158.4667
$ws.Range("J2").Value = 1156
$ws.Range("K2").Value = 1158.4667
$ws.Range("L2").Value = 1156
$ws.Range("M2").Value = -1045.4667
$ws.Range("N2").Value = -1382
$ws.Range("H32").Value = 6066.8184
$ws.Range("I32").Value = 6388.5312
$ws.Range("K32").Value = 6388.5312
$ws.Range("M32").Value = -6101.5312
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""
$ws.Range("H110").Value = 1012.5
$ws.Range("I110").Value = 1037.238
$ws.Range("K110").Value = 1037.238
$ws.Range("M110").Value = 1007.762
$ws.Range("H116").Value = 1157.85
$ws.Range("I116").Value = 1158.4667
$ws.Range("J116").Value = 1156
$ws.Range("K116").Value = 1158.4667
$ws.Range("L116").Value = 1156
$ws.Range("M116").Value = 1135.5333
$ws.Range("N116").Value = -5744

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1157.85
$ws.Range("I3").Value = 1158.4667
$ws.Range("J3").Value = 1156
$ws.Range("K3").Value = 1158.4667
$ws.Range("L3").Value = 1156
$ws.Range("M3").Value = -1044.4667
$ws.Range("N3").Value = -1384
$ws.Range("H86").Value = 1995.7333
$ws.Range("I86").Value = 1813.8235
$ws.Range("K86").Value = 1813.8235
$ws.Range("M86").Value = -690.8235
$ws.Range("H89").Value = 1995.7333
$ws.Range("I89").Value = 1813.8235
$ws.Range("K89").Value = 9069.1175
$ws.Range("M89").Value = -3453.1175
$ws.Range("H94").Value = 883.13336
$ws.Range("I94").Value = 965.2727
$ws.Range("J94").Value = 657.25
$ws.Range("K94").Value = 965.2727
$ws.Range("L94").Value = 657.25
$ws.Range("M94").Value = -514.2727
$ws.Range("N94").Value = -1559.25
$ws.Range("H107").Value = 1162.6923
$ws.Range("I107").Value = 1176.4166
$ws.Range("J107").Value = 998
$ws.Range("K107").Value = 1176.4166
$ws.Range("L107").Value = 998
$ws.Range("M107").Value = 743.5834
$ws.Range("N107").Value = -4838
$ws.Range("H134").Value = 2504.8125
$ws.Range("I134").Value = 1694.16
$ws.Range("J134").Value = 5400
$ws.Range("K134").Value = 5082.48
$ws.Range("L134").Value = 16200
$ws.Range("M134").Value = -2547.48
$ws.Range("N134").Value = -21270

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1607876.2
$ws.Range("I4").Value = 4018452.2
$ws.Range("J4").Value = 825.55554
$ws.Range("K4").Value = 12055356.6
$ws.Range("L4").Value = 2476.66662
$ws.Range("M4").Value = -12055244.6
$ws.Range("N4").Value = -2700.66662
$ws.Range("H6").Value = 83.666664
$ws.Range("I6").Value = 83.666664
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 250.999992
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -137.999992
$ws.Range("N6").Value = ""
$ws.Range("H109").Value = 2200.4614
$ws.Range("J109").Value = 2614.1428
$ws.Range("L109").Value = 7842.428400000001
$ws.Range("N109").Value = -9922.428400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 39800
$ws.Range("J74").Value = 39800
$ws.Range("L74").Value = 39800
$ws.Range("N74").Value = -41672
$ws.Range("H75").Value = 35573.332
$ws.Range("J75").Value = 35573.332
$ws.Range("L75").Value = 35573.332
$ws.Range("N75").Value = -37321.332
$ws.Range("H77").Value = 39800
$ws.Range("J77").Value = 39800
$ws.Range("L77").Value = 119400
$ws.Range("N77").Value = -128760
$ws.Range("H78").Value = 35573.332
$ws.Range("J78").Value = 35573.332
$ws.Range("L78").Value = 106719.996
$ws.Range("N78").Value = -115455.996
$ws.Range("H126").Value = 5430.2324
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5430.2324
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 16290.6972
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = -21230.6972

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5906.25
$ws.Range("I132").Value = 1843.9166
$ws.Range("J132").Value = 11999.75
$ws.Range("K132").Value = 5531.7498
$ws.Range("L132").Value = 35999.25
$ws.Range("M132").Value = -3001.7498
$ws.Range("N132").Value = -41059.25
$ws.Range("H133").Value = 34000
$ws.Range("J133").Value = 34000
$ws.Range("L133").Value = 34000
$ws.Range("N133").Value = -39060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 19000
$ws.Range("J76").Value = 19000
$ws.Range("L76").Value = 19000
$ws.Range("N76").Value = -19630
$ws.Range("H79").Value = 19000
$ws.Range("J79").Value = 19000
$ws.Range("L79").Value = 19000
$ws.Range("N79").Value = -21184
$ws.Range("H113").Value = 501
$ws.Range("I113").Value = 302
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 906
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 1264
$ws.Range("N113").Value = -6440
$ws.Range("H126").Value = 2419.8096
$ws.Range("I126").Value = 1622.7693
$ws.Range("J126").Value = 3715
$ws.Range("K126").Value = 4868.3079
$ws.Range("L126").Value = 11145
$ws.Range("M126").Value = -2398.3079
$ws.Range("N126").Value = -16085
$ws.Range("H136").Value = 7163.2856
$ws.Range("I136").Value = 4286
$ws.Range("J136").Value = 10999.667
$ws.Range("K136").Value = 12858
$ws.Range("L136").Value = 32999.001
$ws.Range("M136").Value = -10308
$ws.Range("N136").Value = -38099.001

Write-Host "Updated 168 cells across sheets"